$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 127.2231523333333
$ws.Range("H2").Value = 381.669457
$ws.Range("I2").Value = 0.6991094561571135
$ws.Range("J2").Value = 0.6991094561571135
$ws.Range("M2").Value = 3.113887
$ws.Range("N2").Value = 9.341661
$ws.Range("O2").Value = 0.8053155340245229
$ws.Range("P2").Value = 0.8053155340245229
$ws.Range("Q2").Value = 396.1585201497864
$ws.Range("R2").Value = 3565.426681348077
$ws.Range("S2").Value = 0.5630037050267597
$ws.Range("T2").Value = 0.5630037050267597

# Row 3
$ws.Range("G3").Value = 127.2231523333333
$ws.Range("H3").Value = 381.669457
$ws.Range("I3").Value = 0.6991094561571135
$ws.Range("J3").Value = 0.6991094561571135
$ws.Range("O3").Value = 0.1946844659754771
$ws.Range("P3").Value = 0.1946844659754771
$ws.Range("Q3").Value = 95.77104461348668
$ws.Range("R3").Value = 861.9394015213801
$ws.Range("S3").Value = 0.1361057511303538
$ws.Range("T3").Value = 0.1361057511303538

# Row 4
$ws.Range("I4").Value = 0.06054748021347716
$ws.Range("J4").Value = 0.06054748021347715
$ws.Range("M4").Value = 3.113887
$ws.Range("N4").Value = 9.341661
$ws.Range("O4").Value = 0.8053155340245229
$ws.Range("P4").Value = 0.8053155340245229
$ws.Range("Q4").Value = 34.30993523105634
$ws.Range("R4").Value = 308.789417079507
$ws.Range("S4").Value = 0.0487598263619556
$ws.Range("T4").Value = 0.04875982636195559

# Row 5
$ws.Range("I5").Value = 0.06054748021347716
$ws.Range("J5").Value = 0.06054748021347715
$ws.Range("O5").Value = 0.1946844659754771
$ws.Range("P5").Value = 0.1946844659754771
$ws.Range("S5").Value = 0.01178765385152157
$ws.Range("T5").Value = 0.01178765385152156

# Row 6
$ws.Range("I6").Value = 0.2403430636294094
$ws.Range("J6").Value = 0.2403430636294094
$ws.Range("M6").Value = 3.113887
$ws.Range("N6").Value = 9.341661
$ws.Range("O6").Value = 0.8053155340245229
$ws.Range("P6").Value = 0.8053155340245229
$ws.Range("Q6").Value = 136.1931977562823
$ws.Range("R6").Value = 1225.738779806541
$ws.Range("S6").Value = 0.1935520026358077
$ws.Range("T6").Value = 0.1935520026358077

# Row 7
$ws.Range("I7").Value = 0.2403430636294094
$ws.Range("J7").Value = 0.2403430636294094
$ws.Range("O7").Value = 0.1946844659754771
$ws.Range("P7").Value = 0.1946844659754771
$ws.Range("S7").Value = 0.04679106099360166
$ws.Range("T7").Value = 0.04679106099360166
